$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ImportProduct")
$ws.Activate() | Out-Null

# Remove the "Giá thuê" column (I) -- "Giá thật" (H) becomes the sole price column, renamed below.
$ws.Columns.Item(9).Delete() | Out-Null

# Remove the "Cho thuê" column (now shifted to P after the delete above).
$ws.Columns.Item(16).Delete() | Out-Null

# Rename remaining headers to match the new "manage image folder" layout.
$ws.Range("B3").Value = "Thương hiệu"
$ws.Range("H3").Value = "Giá bán"
$ws.Range("P3").Value = "Link ảnh chính"
$ws.Range("Q3").Value = "Link ảnh 1"
$ws.Range("R3").Value = "Link ảnh 2"
$ws.Range("S3").Value = "Link ảnh 3"
$ws.Range("T3").Value = "Link ảnh 4"
$ws.Range("U3").Value = "Link ảnh 5"

# Match the new column widths for the newly-visible STT/Thương hiệu columns and the widened image-link column.
$ws.Columns.Item(1).ColumnWidth = 4.5
$ws.Columns.Item(2).ColumnWidth = 11.333333333333334
$ws.Columns.Item(16).ColumnWidth = 15.333333333333334

# Selection moves off the old "select all rows" range onto B3.
$ws.Range("B3").Select() | Out-Null

Write-Output "done"
